$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    "2021-10-05 13:40:54.658435",
    "2021-10-05 13:40:54.658445",
    "2021-10-05 13:40:54.658448",
    "2021-10-05 13:40:54.658450",
    "2021-10-05 13:40:54.658453",
    "2021-10-05 13:40:54.658455",
    "2021-10-05 13:40:54.658458",
    "2021-10-05 13:40:54.658460",
    "2021-10-05 13:40:54.658462",
    "2021-10-05 13:40:54.658464",
    "2021-10-05 13:40:54.658466",
    "2021-10-05 13:40:54.658469",
    "2021-10-05 13:40:54.658471",
    "2021-10-05 13:40:54.658473",
    "2021-10-05 13:40:54.658475",
    "2021-10-05 13:40:54.658477",
    "2021-10-05 13:40:54.658479",
    "2021-10-05 13:40:54.658481",
    "2021-10-05 13:40:54.658483",
    "2021-10-05 13:40:54.658485",
    "2021-10-05 13:40:54.658488",
    "2021-10-05 13:40:54.658489",
    "2021-10-05 13:40:54.658492",
    "2021-10-05 13:40:54.658494",
    "2021-10-05 13:40:54.658496",
    "2021-10-05 13:40:54.658498",
    "2021-10-05 13:40:54.658500",
    "2021-10-05 13:40:54.658502",
    "2021-10-05 13:40:54.658504",
    "2021-10-05 13:40:54.658507",
    "2021-10-05 13:40:54.658509",
    "2021-10-05 13:40:54.658511",
    "2021-10-05 13:40:54.658513",
    "2021-10-05 13:40:54.658516",
    "2021-10-05 13:40:54.658518",
    "2021-10-05 13:40:54.658520",
    "2021-10-05 13:40:54.658522",
    "2021-10-05 13:40:54.658524",
    "2021-10-05 13:40:54.658526",
    "2021-10-05 13:40:54.658529",
    "2021-10-05 13:40:54.658531",
    "2021-10-05 13:40:54.658533",
    "2021-10-05 13:40:54.658535",
    "2021-10-05 13:40:54.658537",
    "2021-10-05 13:40:54.658539",
    "2021-10-05 13:40:54.658542",
    "2021-10-05 13:40:54.658544",
    "2021-10-05 13:40:54.658546",
    "2021-10-05 13:40:54.658548",
    "2021-10-05 13:40:54.658550",
    "2021-10-05 13:40:54.658552",
    "2021-10-05 13:40:54.658554",
    "2021-10-05 13:40:54.658557",
    "2021-10-05 13:40:54.658559",
    "2021-10-05 13:40:54.658561"
)

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value = "time_taken"

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
